$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet from "Sheet" to "Sheet1" ---
$ws.Name = "Sheet1"

# --- New header cells E1:AA1 (A1:D1 already hold Model/Examples/OtherMetrics.../Time (s)) ---
$headers = @(
    "Model", "Examples", "OtherMetrics...", "Time (s)", "Date", "Model Name",
    "Exact Precision (Micro Avg)", "Exact Recall (Micro Avg)", "Exact F1 Score (Micro Avg)",
    "Exact Precision (Macro Avg)", "Exact Recall (Macro Avg)", "Exact F1 Score (Macro Avg)",
    "Exact Precision (Weighted Avg)", "Exact Recall (Weighted Avg)", "Exact F1 Score (Weighted Avg)",
    "Partial Precision", "Partial Recall", "Partial F1 Score",
    "Partial TP", "Partial FP", "Partial FN",
    "Support", "Accuracy", "Result Link", "Stats Link", "No of GPU Used", "Power Consumption"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Bold, boxed (thin border), centered + top-aligned header style for A1:AA1 ---
$headerRange = $ws.Range("A1:AA1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# --- New data row 2 ---
# E2 is a literal text date string, not a real date value -> force text format
# so Excel/the engine doesn't coerce it into a date serial number.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "09/11/2025"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").Value = "Llama-3.3-70B-Instruct"

$ws.Range("G2").Value = 0.3743589743589744
$ws.Range("H2").Value = 0.2723880597014925
$ws.Range("I2").Value = 0.3153347732181426
$ws.Range("J2").Value = 0.1905084681400471
$ws.Range("K2").Value = 0.125914881377906
$ws.Range("L2").Value = 0.1498134368028032
$ws.Range("M2").Value = 0.4157478301135018
$ws.Range("N2").Value = 0.2723880597014925
$ws.Range("O2").Value = 0.3263916762091395
$ws.Range("P2").Value = 0.4845360824742268
$ws.Range("Q2").Value = 0.352059925093633
$ws.Range("R2").Value = 0.4078091106290673

$ws.Range("S2").Value = 94
$ws.Range("T2").Value = 100
$ws.Range("U2").Value = 173
$ws.Range("V2").Value = 268

$ws.Range("W2").Value = 0.9457131136332558

$ws.Range("X2").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Llama-3.3-70B-Instruct_3_shot.txt"
$ws.Range("Y2").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Llama-3.3-70B-Instruct_3_shot.txt"

$ws.Range("Z2").Value = "4 MLGPU"
$ws.Range("AA2").Value = "0.127 kWh"
